$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2305.625
$ws.Range("J17").Value = 2420.7144
$ws.Range("N17").Value = -7598.1432
$ws.Range("L17").Value = 7262.1432
$ws.Range("J32").Value = 9396.357
$ws.Range("M32").Value = -6904.5
$ws.Range("K32").Value = 7230.5
$ws.Range("N32").Value = -10048.357
$ws.Range("L32").Value = 9396.357
$ws.Range("I32").Value = 7230.5
$ws.Range("H32").Value = 8915.056
$ws.Range("I115").Value = 493
$ws.Range("H115").Value = 594.75
$ws.Range("M115").Value = 88
$ws.Range("K115").Value = 1479
$ws.Range("I116").Value = 5001.3335
$ws.Range("H116").Value = 4999.8
$ws.Range("M116").Value = -1559.3335
$ws.Range("K116").Value = 5001.3335
$ws.Range("J138").Value = 0
$ws.Range("M138").Value = -7449.5
$ws.Range("K138").Value = 12589.5
$ws.Range("N138").ClearContents()
$ws.Range("L138").Value = 0
$ws.Range("I138").Value = 4196.5
$ws.Range("H138").Value = 4196.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 1855.375
$ws.Range("H2").Value = 1815.8889
$ws.Range("M2").Value = -1742.375
$ws.Range("K2").Value = 1855.375
$ws.Range("J32").Value = 40000
$ws.Range("N32").Value = -40574
$ws.Range("L32").Value = 40000
$ws.Range("H32").Value = 40000
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("L51").Value = 0
$ws.Range("I116").Value = 1855.375
$ws.Range("H116").Value = 1815.8889
$ws.Range("M116").Value = 438.625
$ws.Range("K116").Value = 1855.375
$ws.Range("I132").Value = 4617.3335
$ws.Range("H132").Value = 6462.625
$ws.Range("M132").Value = -11322.0005
$ws.Range("K132").Value = 13852.0005
$ws.Range("J140").Value = 62000
$ws.Range("N140").Value = -72360
$ws.Range("L140").Value = 62000
$ws.Range("H140").Value = 62000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1815.8889
$ws.Range("M3").Value = -1741.375
$ws.Range("I3").Value = 1855.375
$ws.Range("K3").Value = 1855.375
$ws.Range("I86").Value = 2586.6667
$ws.Range("H86").Value = 4849.75
$ws.Range("M86").Value = -1463.6667
$ws.Range("K86").Value = 2586.6667
$ws.Range("M89").Value = -7317.333500000001
$ws.Range("K89").Value = 12933.3335
$ws.Range("I89").Value = 2586.6667
$ws.Range("H89").Value = 4849.75
$ws.Range("I107").Value = 8039.1333
$ws.Range("L107").Value = 3500
$ws.Range("H107").Value = 7505.1177
$ws.Range("J107").Value = 3500
$ws.Range("M107").Value = -6119.1333
$ws.Range("K107").Value = 8039.1333
$ws.Range("N107").Value = -7340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M16").Value = -716.8333
$ws.Range("K16").Value = 1003.8333
$ws.Range("I16").Value = 1003.8333
$ws.Range("H16").Value = 1094.5555
$ws.Range("I31").Value = 2798.5
$ws.Range("H31").Value = 61398.25
$ws.Range("M31").Value = -2503.5
$ws.Range("K31").Value = 2798.5
$ws.Range("I34").Value = 2798.5
$ws.Range("H34").Value = 61398.25
$ws.Range("M34").Value = -2596.5
$ws.Range("K34").Value = 2798.5
$ws.Range("K58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("H58").Value = 7500
$ws.Range("M58").ClearContents()
$ws.Range("I107").Value = 689.2857
$ws.Range("H107").Value = 748.5
$ws.Range("M107").Value = 1230.7143
$ws.Range("K107").Value = 689.2857
$ws.Range("I113").Value = 1003.8333
$ws.Range("H113").Value = 1094.5555
$ws.Range("M113").Value = 1166.1667
$ws.Range("K113").Value = 1003.8333
$ws.Range("J122").Value = 1510.25
$ws.Range("N122").Value = -9430.75
$ws.Range("L122").Value = 4530.75
$ws.Range("H122").Value = 1699.5555
$ws.Range("K136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("H136").Value = 7500
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K14").Value = 1747.66662
$ws.Range("I14").Value = 582.55554
$ws.Range("H14").Value = 582.55554
$ws.Range("M14").Value = -1574.66662
$ws.Range("I44").Value = 87.5
$ws.Range("H44").Value = 87.5
$ws.Range("M44").Value = 135.5
$ws.Range("K44").Value = 262.5
$ws.Range("J109").Value = 0
$ws.Range("M109").Value = -1594
$ws.Range("K109").Value = 2634
$ws.Range("N109").ClearContents()
$ws.Range("L109").Value = 0
$ws.Range("I109").Value = 878
$ws.Range("H109").Value = 878
$ws.Range("M131").Value = -15739.9995
$ws.Range("K131").Value = 20779.9995
$ws.Range("I131").Value = 6926.6665
$ws.Range("H131").Value = 3567.7144
$ws.Range("K139").Value = 7245
$ws.Range("N139").Value = -22280
$ws.Range("L139").Value = 12000
$ws.Range("I139").Value = 2415
$ws.Range("H139").Value = 2943.3333
$ws.Range("J139").Value = 4000
$ws.Range("M139").Value = -2105

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 181503330
$ws.Range("M11").Value = -181503191
$ws.Range("K11").Value = 181503330
$ws.Range("I11").Value = 181503330
$ws.Range("I107").Value = 170.14285
$ws.Range("H107").Value = 386.22223
$ws.Range("M107").Value = 1749.85715
$ws.Range("K107").Value = 170.14285
$ws.Range("N113").ClearContents()
$ws.Range("I113").Value = 1066.6666
$ws.Range("L113").Value = 0
$ws.Range("H113").Value = 1066.6666
$ws.Range("J113").Value = 0
$ws.Range("M113").Value = 1103.3334
$ws.Range("K113").Value = 1066.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 527.6
$ws.Range("H22").Value = 639.6667
$ws.Range("M22").Value = -232.6
$ws.Range("K22").Value = 527.6
$ws.Range("H27").Value = 639.6667
$ws.Range("M27").Value = -420.6
$ws.Range("K27").Value = 527.6
$ws.Range("I27").Value = 527.6
$ws.Range("L33").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("L46").Value = 2065.1177
$ws.Range("I46").Value = 1026.3125
$ws.Range("H46").Value = 1561.4546
$ws.Range("J46").Value = 2065.1177
$ws.Range("M46").Value = -838.3125
$ws.Range("K46").Value = 1026.3125
$ws.Range("N46").Value = -2441.1177
$ws.Range("J93").Value = 200
$ws.Range("M93").Value = -336.7141999999999
$ws.Range("K93").Value = 1584.7142
$ws.Range("N93").Value = -2696
$ws.Range("L93").Value = 200
$ws.Range("I93").Value = 1584.7142
$ws.Range("H93").Value = 1411.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K136").Value = 8143.600199999999
$ws.Range("I136").Value = 2714.5334
$ws.Range("H136").Value = 2701.125
$ws.Range("M136").Value = -5593.600199999999
